# Add 2022-Q3 data.
#
# Before: sheet1 "总计" (summary), sheet2 "2022-Q2" (per-fund detail for Q2).
# After:  sheet1 "总计" (summary, now with a Q3 row and the old Q2 row),
#         sheet2 "2022-Q3" (per-fund detail for Q3, reusing the old Q2 slot),
#         sheet3 "2022-Q2" (per-fund detail for Q2, moved to a new sheet).

function Set-TextValue($range, $value) {
    # Force a numeric-looking string (e.g. a fund code like "006263", or a
    # decimal string like "0.20") to be stored as text instead of being
    # auto-converted to a number, without leaving any numFmt/style residue
    # behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# 1. Duplicate the existing "2022-Q2" detail sheet so the old data keeps
#    living on its own tab (this becomes the new sheet3 / rId3).
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)

# 2. Turn the original "2022-Q2" sheet into the "2022-Q3" sheet (rename it
#    first, before renaming the copy back to "2022-Q2", to avoid a name
#    clash) and overwrite its data with the Q3 numbers, restyling the
#    header/A-column cells to match the style used elsewhere for freshly
#    written sheets (same style as the summary sheet's header cells).
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

Set-TextValue $q2.Range("B2") "161124"
$q2.Range("C2").Value = "易方达香港恒生综合小型股指数（QDII-LOF）A"
Set-TextValue $q2.Range("D2") "0.20"
Set-TextValue $q2.Range("E2") "91.61"
Set-TextValue $q2.Range("F2") "1.29"
Set-TextValue $q2.Range("G2") "0.0026"
$q2.Range("H2").Value = 8

Set-TextValue $q2.Range("B3") "006263"
$q2.Range("C3").Value = "易方达香港恒生综合小型股指数（QDII-LOF）C"
Set-TextValue $q2.Range("D3") "0.05"
Set-TextValue $q2.Range("E3") "91.61"
Set-TextValue $q2.Range("F3") "1.29"
Set-TextValue $q2.Range("G3") "0.0006"
$q2.Range("H3").Value = 8

$summary.Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q2.Range("A2:A3").PasteSpecial(-4122)

# 3. Update the summary sheet: row 2 becomes the 2022-Q3 entry, and a new
#    row 3 is added holding the data that used to be in row 2 (2022-Q2).
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.01

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0
